# Refresh the "times" data in Sheet1 with the latest run-time figures
# (mirrors a Data > Refresh All of the times.csv text connection) and
# clear the yellow "needs review" highlight from the rows that were
# updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated SECONDS column values pulled from the refreshed times.csv ---
$updates = @{
    "B19" = 27.5     # MN
    "B23" = 14.4     # NV
    "B25" = 12.1     # NM
    "B26" = 541.7    # NY
    "B30" = 15.2     # OR
    "B33" = 14.6     # TN
    "B37" = 50.1     # WA
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Clear the yellow highlight fill from the state-name cells in column A
#     for the rows whose data just got refreshed ---
$highlightedRows = @(19, 23, 25, 26, 30, 33, 37)
foreach ($r in $highlightedRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Interior.ColorIndex = -4142
    $cell.Interior.Pattern = -4142
}

# --- Move the selection to reflect where the editor left off ---
$null = $ws.Range("H21").Select()
